# #5: cash & deposit done
# Adds structured metadata columns (G:M) to the "存款" (deposits) sheet and
# turns row 1 into a real header row (bank, deposit_type, currency, owner,
# total, property_category, category, date, legislator_name, legislator_id,
# source_file, index), mirroring the property_category/category/date/
# legislator_name/legislator_id/source_file/index fields onto every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Copy the header cell formatting (bold font + border) onto the new cells.
$ws.Range("F1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Force the "date" column to plain text first so that the literal
# "2012-05-01" strings typed below are not auto-converted into date
# serial numbers by Excel's input parser.
$ws.Range("I2:I9").NumberFormat = "@"

# ---- Data rows (2-9) -----------------------------------------------------
$banks = @("中國信託商業銀行斗六分行","復華商業銀行","臺灣中小企業銀行","華南商業銀行斗六分行","合作金庫商業銀行","京城商業銀行","雲林縣斗六市農會","臺灣銀行斗六分行")
$totals = @(49212, 278915, 145000, 3641580, 275, 395736, 303891, 2506288)
$indexVals = @(43, 44, 46, 47, 48, 49, 50, 51)

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $banks[$i]          # B: bank
    $ws.Cells.Item($r, 3).Value = "活期儲蓄存款"        # C: deposit_type
    $ws.Cells.Item($r, 4).Value = "新臺幣"              # D: currency
    $ws.Cells.Item($r, 5).Value = "劉建國"              # E: owner
    $ws.Cells.Item($r, 6).Value = $totals[$i]          # F: total
    $ws.Cells.Item($r, 7).Value = "deposit"            # G: property_category
    $ws.Cells.Item($r, 8).Value = "normal"             # H: category
    $ws.Cells.Item($r, 9).Value = "2012-05-01"         # I: date
    $ws.Cells.Item($r, 10).Value = "劉建國"             # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 1723                # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmpd6c01"          # L: source_file
    $ws.Cells.Item($r, 13).Value = $indexVals[$i]      # M: index
}

# Restore the "date" column back to the same (General) formatting used by
# the rest of the data rows, now that the literal text values are in place.
$ws.Range("F2").Copy()
$ws.Range("I2:I9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wb.Save()
